$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4's data replaces row 3's data (A3 and D3 and I3 change),
# then row 4 is removed entirely.
$ws.Range("A3").Value = "Test"
$ws.Range("D3").Value = 44.73
$ws.Range("I3").Value = 3

$ws.Rows.Item(4).Delete()
